$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 10 (pushes existing rows 10:41 down to 11:42,
# and extends the used range to A1:R42). Excel carries the formatting
# (e.g. the date style on column D) down from the row above automatically.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row with the new weekly price observation.
$ws.Cells.Item(10, 1).Value  = 9
$ws.Cells.Item(10, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(10, 3).Value  = "Metropolitana"
$ws.Cells.Item(10, 4).Value  = 44971
$ws.Cells.Item(10, 5).Value  = 13
$ws.Cells.Item(10, 6).Value  = 100112010
$ws.Cells.Item(10, 7).Value  = "Achicoria"
$ws.Cells.Item(10, 8).Value  = "Sin especificar"
$ws.Cells.Item(10, 9).Value  = "Primera"
$ws.Cells.Item(10, 10).Value = 160
$ws.Cells.Item(10, 11).Value = 7000
$ws.Cells.Item(10, 12).Value = 7000
$ws.Cells.Item(10, 13).Value = 7000
$ws.Cells.Item(10, 14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(10, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(10, 16).Value = 438
$ws.Cells.Item(10, 17).Value = 16
$ws.Cells.Item(10, 18).Value = "Hortaliza"
